# Manual_Scenarios_TestCases-Share: content refresh
# - "Test Case" sheet: update the Steps-to-Take text for TC01 (new wording,
#   mentions navigating to the Product module / thumbnails) and grow row 2 to fit it.
# - "Test Scenario" sheet: extend the red "Without selecting Contacts..." note with
#   more negative-testing bullets (change layout / search / filter while sharing).
# - Make "Test Scenario" the active/selected sheet (was "Test Case").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Test Case" sheet - update Steps to Take (H2) and row height.
# ---------------------------------------------------------------------------
$wsTC = $wb.Worksheets.Item("Test Case")

$newSteps = '1. Log in to "http://blubox.shoppinpal.com/s eller/" website' + "`n" `
  + '2. Navigate to Product module. Verify the following:' + "`n" `
  + 'a. User is able to see thumbnails in the selected layout.' + "`n" `
  + "b. User's E-mail address is displaying correctly" + "`n" `
  + "c. User can't edit these details" + "`n" `
  + 'd. Various tabs are present like: Products, Catalog, Gallery, Contact'

$wsTC.Range("H2").Value = $newSteps

# Row grows to fit the extra wrapped line.
$wsTC.Rows.Item(2).RowHeight = 165

# ---------------------------------------------------------------------------
# 2) "Test Scenario" sheet - extend the second (black) run of the rich-text
#    note in E8 / E9 while keeping the first (red) run untouched.
# ---------------------------------------------------------------------------
$wsTS = $wb.Worksheets.Item("Test Scenario")

$newNote = 'while sharing to contact try to delete contact, while sharing to contact change layout,  while sharing to contact do a search, while sharing to contact put filter,  '

$cellE8 = $wsTS.Range("E8")
$cellE8.Characters(42, 39).Text = $newNote
$cellE8.Characters(1, 41).Font.Color = 255
$cellE8.Characters(42, $newNote.Length).Font.Color = 0

# E9 carries the exact same note as E8 - copy the formatted rich text across
# so both cells keep sharing a single string entry (as in the original file).
$cellE8.Copy()
$wsTS.Range("E9").PasteSpecial()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Make "Test Scenario" the active sheet / selection (was "Test Case").
# ---------------------------------------------------------------------------
$wsTS.Activate()
$wsTS.Range("E8").Select()
